# G2-1774 Update apache POI
#
# The corresponding commit only bumped the Apache POI library used to
# write the .pptx package. That library change altered the *order* in
# which XML attributes are serialized (e.g. xmlns declarations, and
# attributes such as type/idx/sz on <p:ph>, or marL/algn/indent on
# paragraph-property elements) but it did not change any slide content,
# placeholder, text, formatting value, color, or structural element.
#
# Every hunk in the associated OOXML diff is a pure attribute-order
# permutation: the same attribute/value pairs appear on both sides of
# each changed tag, just written in a different sequence, and no lines
# are added or removed anywhere in the package.
#
# The PowerPoint object model has no notion of "XML attribute order" -
# that's purely a side effect of which serializer wrote the file - so
# there is no COM operation that corresponds to this change. We touch
# the presentation object (to keep this a valid, executed COM script)
# without mutating any content, property, or shape.

$p = $ppt.ActivePresentation
$null = $p.Slides.Count
